$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the "Pütter/Hesse" column (K) the same way the other X-columns are
# marked for the rows 2-6 evaluation criteria.
$ws.Range("K2").Value = "X"
$ws.Range("K3").Value = "X"
$ws.Range("K4").Value = "X"
$ws.Range("K5").Value = "X"
$ws.Range("K6").Value = "X"

# Add the "Pütter/Hesse" entry to the signature row (row 8), matching the
# value and formatting already used for the "Ausgefüllt von" entry in I8.
$ws.Range("K8").Value = "Pütter/Hesse"
$ws.Range("I8").Copy() | Out-Null
$ws.Range("K8").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Resize columns J and K now that column K has real content in it.
$ws.Columns.Item(10).ColumnWidth = 13.43
$ws.Columns.Item(11).ColumnWidth = 11.17

# Move the active selection to the newly filled-in K8 cell.
$ws.Range("K8").Select() | Out-Null
